$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 535.3333
$ws.Range("I28").Value = 629.04346
$ws.Range("J28").Value = 227.42857
$ws.Range("K28").Value = 629.04346
$ws.Range("L28").Value = 227.42857
$ws.Range("M28").Value = -144.04346
$ws.Range("N28").Value = -1197.42857

$ws.Range("H96").Value = 283.33334
$ws.Range("I96").Value = 275
$ws.Range("J96").Value = 300
$ws.Range("K96").Value = 825
$ws.Range("L96").Value = 900
$ws.Range("M96").Value = 548
$ws.Range("N96").Value = -3646

$ws.Range("H98").Value = 1707.5
$ws.Range("I98").Value = 1075.4166
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 1075.4166
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = 422.5834

$ws.Range("H122").Value = 1707.5
$ws.Range("I122").Value = 1075.4166
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 3226.2498
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -776.2498000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2271.5715
$ws.Range("I2").Value = 1414.7727
$ws.Range("J2").Value = 5413.1665
$ws.Range("K2").Value = 1414.7727
$ws.Range("L2").Value = 5413.1665
$ws.Range("M2").Value = -1301.7727

$ws.Range("H7").Value = 30622.223
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 30622.223
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 30622.223
$ws.Range("N7").Value = -30850.223

$ws.Range("H32").Value = 2705263.5
$ws.Range("I32").Value = 3702.0632
$ws.Range("J32").Value = 17949790
$ws.Range("K32").Value = 3702.0632
$ws.Range("L32").Value = 17949790
$ws.Range("M32").Value = -3415.0632
$ws.Range("N32").Value = -17950364

$ws.Range("H45").Value = 3159.9614
$ws.Range("I45").Value = 2342.6
$ws.Range("J45").Value = 4274.5454
$ws.Range("K45").Value = 2342.6
$ws.Range("L45").Value = 4274.5454
$ws.Range("M45").Value = -1965.6
$ws.Range("N45").Value = -5028.5454

$ws.Range("H52").Value = 16741.54
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 16741.54
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 16741.54
$ws.Range("N52").Value = -17377.54

$ws.Range("H110").Value = 10321.772
$ws.Range("I110").Value = 10559.223
$ws.Range("J110").Value = 9253.25
$ws.Range("K110").Value = 10559.223
$ws.Range("L110").Value = 9253.25
$ws.Range("M110").Value = -8514.223
$ws.Range("N110").Value = -13343.25

$ws.Range("H116").Value = 2271.5715
$ws.Range("I116").Value = 1414.7727
$ws.Range("J116").Value = 5413.1665
$ws.Range("K116").Value = 1414.7727
$ws.Range("L116").Value = 5413.1665
$ws.Range("M116").Value = 879.2273

$ws.Range("H132").Value = 110001.22
$ws.Range("I132").Value = 120189.48
$ws.Range("J132").Value = 3024.5
$ws.Range("K132").Value = 360568.44
$ws.Range("L132").Value = 9073.5
$ws.Range("M132").Value = -358038.44
$ws.Range("N132").Value = -14133.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2271.5715
$ws.Range("I3").Value = 1414.7727
$ws.Range("J3").Value = 5413.1665
$ws.Range("K3").Value = 1414.7727
$ws.Range("L3").Value = 5413.1665
$ws.Range("M3").Value = -1300.7727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1750.25
$ws.Range("I22").Value = 2000.3334
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2000.3334
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1650.3334
$ws.Range("N22").Value = -1700

$ws.Range("H31").Value = 1654.25
$ws.Range("I31").Value = 1222.5
$ws.Range("J31").Value = 3122.2
$ws.Range("K31").Value = 1222.5
$ws.Range("L31").Value = 3122.2
$ws.Range("M31").Value = -927.5
$ws.Range("N31").Value = -3712.2

$ws.Range("H34").Value = 1654.25
$ws.Range("I34").Value = 1222.5
$ws.Range("J34").Value = 3122.2
$ws.Range("K34").Value = 1222.5
$ws.Range("L34").Value = 3122.2
$ws.Range("M34").Value = -1020.5
$ws.Range("N34").Value = -3526.2

$ws.Range("H64").Value = 33000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 33000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 33000
$ws.Range("N64").Value = -33496
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 33000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 33000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 33000
$ws.Range("N67").Value = -34716
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 696.5227
$ws.Range("I113").Value = 594.75
$ws.Range("J113").Value = 734.6875
$ws.Range("K113").Value = 1784.25
$ws.Range("L113").Value = 2204.0625
$ws.Range("M113").Value = 385.75
$ws.Range("N113").Value = -6544.0625

$ws.Range("H122").Value = 47611.71
$ws.Range("I122").Value = 393.75
$ws.Range("J122").Value = 52218.34
$ws.Range("K122").Value = 3543.75
$ws.Range("L122").Value = 469965.0599999999
$ws.Range("M122").Value = -1093.75
$ws.Range("N122").Value = -474865.0599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 9020.236999999999
$ws.Range("I126").Value = 2694.5
$ws.Range("J126").Value = 14713.4
$ws.Range("K126").Value = 8083.5
$ws.Range("L126").Value = 44140.2
$ws.Range("M126").Value = -5613.5
$ws.Range("N126").Value = -49080.2

$ws.Range("H132").Value = 2102.238
$ws.Range("I132").Value = 1681.9131
$ws.Range("J132").Value = 2611.0527
$ws.Range("K132").Value = 5045.7393
$ws.Range("L132").Value = 7833.158100000001
$ws.Range("M132").Value = -2515.7393
$ws.Range("N132").Value = -12893.1581

$ws.Range("H136").Value = 20442.084
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 20442.084
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 61326.25199999999
$ws.Range("N136").Value = -66426.25199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3195.7144
$ws.Range("I122").Value = 2411.111
$ws.Range("J122").Value = 3784.1667
$ws.Range("K122").Value = 7233.333
$ws.Range("L122").Value = 11352.5001
$ws.Range("M122").Value = -4783.333
$ws.Range("N122").Value = -16252.5001

$ws.Range("H132").Value = 2727.9688
$ws.Range("I132").Value = 2220.5
$ws.Range("J132").Value = 4250.375
$ws.Range("K132").Value = 6661.5
$ws.Range("L132").Value = 12751.125
$ws.Range("M132").Value = -4131.5

$ws.Range("H133").Value = 39659.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39659.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39659.332
$ws.Range("N133").Value = -44719.332

$ws.Range("H136").Value = 21864.22
$ws.Range("I136").Value = 1673.9524
$ws.Range("J136").Value = 127863.125
$ws.Range("K136").Value = 5021.857199999999
$ws.Range("L136").Value = 383589.375
$ws.Range("M136").Value = -2471.857199999999
$ws.Range("N136").Value = -388689.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 29500
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 29500
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 29500
$ws.Range("N119").Value = -39176

$ws.Range("H122").Value = 3371.25
$ws.Range("I122").Value = 2000.6
$ws.Range("J122").Value = 4741.9
$ws.Range("K122").Value = 6001.799999999999
$ws.Range("L122").Value = 14225.7
$ws.Range("M122").Value = -3551.799999999999

$ws.Range("H132").Value = 1556.4894
$ws.Range("I132").Value = 1389.9318
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 4169.7954
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -1639.7954
$ws.Range("N132").Value = -17057.9999

$ws.Range("H133").Value = 39191.668
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39191.668
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39191.668
$ws.Range("N133").Value = -49311.668
